$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update detected-location coordinates in column B per the commit diff
$ws.Range("B7").Value = "32.766882,34.967053"
$ws.Range("B8").Value = "32.511729,35.502029"
$ws.Range("B10").Value = "32.267628,34.993511"
$ws.Range("B11").Value = "32.792761,34.995336"
$ws.Range("B18").Value = "32.980490,35.542420"
$ws.Range("B19").Value = "31.960770,34.876512"
$ws.Range("B21").Value = "32.986934,35.708518"
$ws.Range("B24").Value = "33.005860,35.094090"
$ws.Range("B25").Value = "32.601426,35.289751"
$ws.Range("B29").Value = "33.194459,35.572940"
$ws.Range("B30").Value = "31.961063,34.807761"
$ws.Range("B31").Value = "32.036425,34.842884"
$ws.Range("B32").Value = "31.977527,34.808252"
$ws.Range("B34").Value = "31.784215,35.117210"
$ws.Range("B35").Value = "32.058998,34.815227"
$ws.Range("B36").Value = "32.121447,34.803699"
$ws.Range("B37").Value = "32.692764,34.940222"
$ws.Range("B41").Value = "32.775683,34.967878"
$ws.Range("B42").Value = "32.199671,35.212911"
$ws.Range("B43").Value = "32.049533,34.764483"
$ws.Range("B51").Value = "32.792761,34.995336"
$ws.Range("B56").Value = "31.814560,34.779980"
$ws.Range("B57").Value = "31.248833,35.198232"
$ws.Range("B58").Value = "32.074578,34.805974"
$ws.Range("B61").Value = "32.590574,34.936472"
$ws.Range("B64").Value = "32.018460,34.748167"
$ws.Range("B69").Value = "31.822668,35.253867"
$ws.Range("B70").Value = "31.791658,34.651074"
$ws.Range("B75").Value = "31.916670,35.016670"
$ws.Range("B77").Value = "32.095980,34.774333"
$ws.Range("B83").Value = "31.248833,35.198232"
$ws.Range("B85").Value = "31.945204,34.878075"
$ws.Range("B88").Value = "32.098181,34.896471"
$ws.Range("B93").Value = "31.750898,35.207819"
$ws.Range("B96").Value = "32.177911,34.905656"
$ws.Range("B97").Value = "32.049272,34.798714"
$ws.Range("B104").Value = "32.009918,34.739188"
$ws.Range("B109").Value = "31.314100,34.620250"
$ws.Range("B111").Value = "32.860863,35.099385"
$ws.Range("B114").Value = "31.977527,34.808252"
$ws.Range("B116").Value = "31.682230,34.745240"
$ws.Range("B119").Value = "31.670900,34.779750"
$ws.Range("B120").Value = "32.093309,34.885509"
$ws.Range("B141").Value = "32.447787,34.914428"
$ws.Range("B145").Value = "32.139558,34.959151"
$ws.Range("B168").Value = "32.696543,35.052512"
$ws.Range("B169").Value = "31.746624,35.213931"
$ws.Range("B176").Value = "32.819327,34.997607"
$ws.Range("B183").Value = "32.603829,35.298516"
$ws.Range("B184").Value = "32.807619,35.057422"
$ws.Range("B185").Value = "32.804853,35.072857"
$ws.Range("B187").Value = "32.607430,35.292122"
$ws.Range("B188").Value = "32.782155,34.976622"
$ws.Range("B189").Value = "32.054678,34.804617"
$ws.Range("B195").Value = "32.815143,35.060720"
$ws.Range("B199").Value = "32.471755,34.969755"
$ws.Range("B205").Value = "32.163217,34.961133"
$ws.Range("B206").Value = "32.163217,34.961133"
$ws.Range("B209").Value = "32.790077,35.516279"
$ws.Range("B211").Value = "32.139558,34.959151"
$ws.Range("B213").Value = "31.785087,35.210391"
$ws.Range("B224").Value = "32.139558,34.959151"
$ws.Range("B243").Value = "32.916364,35.125162"
$ws.Range("B247").Value = "32.777129,35.040632"
$ws.Range("B254").Value = "32.854815,35.196680"
$ws.Range("B255").Value = "32.163217,34.961133"
$ws.Range("B258").Value = "31.750585,35.215673"
$ws.Range("B267").Value = "31.750859,35.213920"
$ws.Range("B268").Value = "32.139558,34.959151"
$ws.Range("B270").Value = "32.099592,34.828768"
$ws.Range("B275").Value = "32.046554,34.869660"
$ws.Range("B276").Value = "32.984064,35.248787"
$ws.Range("B282").Value = "31.813664,34.667678"
$ws.Range("B283").Value = "31.663407,34.599960"
$ws.Range("B289").Value = "32.818629,34.996670"
$ws.Range("B295").Value = "31.750492,35.215772"
$ws.Range("B296").Value = "31.785775,35.182708"
$ws.Range("B297").Value = "31.750585,35.215673"
$ws.Range("B299").Value = "31.750492,35.215772"
$ws.Range("B307").Value = "32.194693,34.884294"
$ws.Range("B308").Value = "31.749963,35.141489"
$ws.Range("B312").Value = "33.003241,35.091790"
$ws.Range("B316").Value = "31.750859,35.213920"
$ws.Range("B335").Value = "32.949509,35.174243"
$ws.Range("B336").Value = "31.784988,35.210374"
$ws.Range("B340").Value = "31.757498,35.218264"
$ws.Range("B349").Value = "31.897055,34.800408"
$ws.Range("B354").Value = "32.163217,34.961133"
$ws.Range("B356").Value = "32.151070,34.847113"
$ws.Range("B358").Value = "31.652206,34.579598"
$ws.Range("B361").Value = "31.747225,35.212499"
$ws.Range("B365").Value = "31.773929,34.629620"
$ws.Range("B373").Value = "31.784215,35.117210"
$ws.Range("B375").Value = "32.151070,34.847113"
$ws.Range("B376").Value = "32.100120,34.828677"
$ws.Range("B384").Value = "32.139558,34.959151"
$ws.Range("B387").Value = "32.139558,34.959151"
$ws.Range("B388").Value = "32.165553,34.813406"
$ws.Range("B390").Value = "31.982527,34.765084"
$ws.Range("B392").Value = "31.807623,34.664804"
$ws.Range("B398").Value = "32.055436,34.805472"
$ws.Range("B400").Value = "32.020682,34.805150"
$ws.Range("B401").Value = "32.097022,34.829235"
$ws.Range("B403").Value = "31.665784,34.601137"
$ws.Range("B407").Value = "31.677567,34.596921"
$ws.Range("B410").Value = "31.946849,34.879864"
$ws.Range("B415").Value = "31.858484,35.215449"
$ws.Range("B417").Value = "31.750492,35.215772"
$ws.Range("B419").Value = "31.858484,35.215449"
$ws.Range("B420").Value = "31.792463,35.144323"
$ws.Range("B421").Value = "32.163217,34.961133"
$ws.Range("B423").Value = "32.163217,34.961133"
$ws.Range("B426").Value = "31.238529,34.795441"
$ws.Range("B429").Value = "32.163217,34.961133"
$ws.Range("B438").Value = "31.319243,34.623672"
$ws.Range("B453").Value = "31.826014,34.658552"
$ws.Range("B469").Value = "32.163217,34.961133"
$ws.Range("B476").Value = "32.139558,34.959151"
$ws.Range("B480").Value = "32.054678,34.804617"
$ws.Range("B483").Value = "32.047811,34.882122"
